# Disaggregation of commodity Copper
# 1. Rename the commodity label "Copper ores and concentrates" -> "Copper"
#    (this text is a shared string used identically in cell C4 on every
#    year sheet, so it must be updated on every sheet so the shared
#    string is fully replaced rather than leaving a duplicate string
#    behind).
# 2. A handful of sheets also carry a tiny floating point re-rounding of
#    the value in D4 (last significant digit only).

$wb = $excel.ActiveWorkbook

$sheetCount = $wb.Worksheets.Count
for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Range("C4").Value = "Copper"
}

# Map of worksheet index -> updated D4 value (precision-only changes)
$d4Updates = @{
    42 = 556760.218623249
    68 = 1496659.553132901
    74 = 3038611.237094068
    76 = 3326325.209350231
    78 = 3483577.984832576
    99 = 3544367.360995423
}

foreach ($idx in $d4Updates.Keys) {
    $ws = $wb.Worksheets.Item([int]$idx)
    $ws.Range("D4").Value = $d4Updates[$idx]
}

Write-Output "Updated Copper label on $sheetCount sheets and $($d4Updates.Count) D4 values."
